$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 9, pushing the existing rows 9-14 down to 11-16.
$ws.Rows("9:10").Insert()

# Row 9: new weekly entry (Primera / $/bandeja 10 kilos)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44907
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100103
$ws.Range("H9").Value = "Frutos de hueso (carozo)"
$ws.Range("I9").Value = 100103003
$ws.Range("J9").Value = "Damasco"
$ws.Range("K9").Value = "Castle Brite"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 15500
$ws.Range("Q9").Value = "`$/bandeja 10 kilos"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 1550
$ws.Range("T9").Value = 10

# Row 10: new weekly entry (Segunda / $/bandeja 10 kilos)
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44907
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100103
$ws.Range("H10").Value = "Frutos de hueso (carozo)"
$ws.Range("I10").Value = 100103003
$ws.Range("J10").Value = "Damasco"
$ws.Range("K10").Value = "Castle Brite"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 14000
$ws.Range("P10").Value = 14000
$ws.Range("Q10").Value = "`$/bandeja 10 kilos"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 1400
$ws.Range("T10").Value = 10
